$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3781.9
$arr[0,1] = 1935
$arr[0,2] = 5628.8
$arr[0,3] = 5805
$arr[0,4] = 16886.4
$arr[0,5] = -5524
$arr[0,6] = -17448.4
$wsALC.Range("H29:N29").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 10000
$arr[0,1] = 0
$arr[0,2] = 10000
$arr[0,3] = 0
$arr[0,4] = 10000
$arr[0,5] = $null
$arr[0,6] = -10496
$wsALC.Range("H64:N64").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 10000
$arr[0,1] = 0
$arr[0,2] = 10000
$arr[0,3] = 0
$arr[0,4] = 10000
$arr[0,5] = $null
$arr[0,6] = -11716
$wsALC.Range("H67:N67").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 26950
$arr[0,1] = 3900
$arr[0,2] = 50000
$arr[0,3] = 3900
$arr[0,4] = 50000
$arr[0,5] = -2964
$arr[0,6] = -51872
$wsALC.Range("H74:N74").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 26950
$arr[0,1] = 3900
$arr[0,2] = 50000
$arr[0,3] = 19500
$arr[0,4] = 250000
$arr[0,5] = -14820
$arr[0,6] = -259360
$wsALC.Range("H77:N77").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 550
$arr[0,1] = 100
$arr[0,2] = 1000
$arr[0,3] = 300
$arr[0,4] = 3000
$arr[0,5] = 698
$arr[0,6] = -4996
$wsALC.Range("H80:N80").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 550
$arr[0,1] = 100
$arr[0,2] = 1000
$arr[0,3] = 900
$arr[0,4] = 9000
$arr[0,5] = 4092
$arr[0,6] = -18984
$wsALC.Range("H83:N83").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 985
$arr[0,1] = 985
$arr[0,2] = 0
$arr[0,3] = 2955
$arr[0,4] = 0
$arr[0,5] = -1388
$arr[0,6] = $null
$wsALC.Range("H115:N115").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1329.3334
$arr[0,1] = 1021.6667
$arr[0,2] = 1944.6666
$arr[0,3] = 3065.0001
$arr[0,4] = 5833.9998
$arr[0,5] = 2074.9999
$arr[0,6] = -16113.9998
$wsALC.Range("H138:N138").Value = $arr

$wsARM = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1238.6
$arr[0,1] = 1104.4706
$arr[0,2] = 1998.6666
$arr[0,3] = 1104.4706
$arr[0,4] = 1998.6666
$arr[0,5] = -817.4706000000001
$arr[0,6] = -2572.6666
$wsARM.Range("H32:N32").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 55000
$arr[0,1] = 0
$arr[0,2] = 55000
$arr[0,3] = 0
$arr[0,4] = 55000
$arr[0,5] = $null
$arr[0,6] = -64960
$wsARM.Range("H128:N128").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 40709
$arr[0,1] = 40709
$arr[0,2] = 0
$arr[0,3] = 40709
$arr[0,4] = 0
$arr[0,5] = -35709
$arr[0,6] = $null
$wsARM.Range("H129:N129").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4410.5
$arr[0,1] = 5489.2
$arr[0,2] = 1713.75
$arr[0,3] = 16467.6
$arr[0,4] = 5141.25
$arr[0,5] = -13937.6
$arr[0,6] = -10201.25
$wsARM.Range("H132:N132").Value = $arr

$wsBSM = $wb.Worksheets.Item("BSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 8447.556
$arr[0,1] = 7916.4546
$arr[0,2] = 9282.143
$arr[0,3] = 7916.4546
$arr[0,4] = 9282.143
$arr[0,5] = -5996.4546
$arr[0,6] = -13122.143
$wsBSM.Range("H107:N107").Value = $arr

$wsCRP = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3917
$arr[0,1] = 1885.8
$arr[0,2] = 8995
$arr[0,3] = 1885.8
$arr[0,4] = 8995
$arr[0,5] = -1682.8
$arr[0,6] = -9401
$wsCRP.Range("H58:N58").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1662.8
$arr[0,1] = 1662.8
$arr[0,2] = 0
$arr[0,3] = 1662.8
$arr[0,4] = 0
$arr[0,5] = -1038.8
$arr[0,6] = $null
$wsCRP.Range("H62:N62").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1662.8
$arr[0,1] = 1662.8
$arr[0,2] = 0
$arr[0,3] = 8314
$arr[0,4] = 0
$arr[0,5] = -5194
$arr[0,6] = $null
$wsCRP.Range("H65:N65").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = $null
$arr[0,6] = $null
$wsCRP.Range("H132:N132").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 110000
$arr[0,1] = 0
$arr[0,2] = 110000
$arr[0,3] = 0
$arr[0,4] = 110000
$arr[0,5] = $null
$arr[0,6] = -120140
$wsCRP.Range("H135:N135").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3917
$arr[0,1] = 1885.8
$arr[0,2] = 8995
$arr[0,3] = 5657.4
$arr[0,4] = 26985
$arr[0,5] = -3107.4
$arr[0,6] = -32085
$wsCRP.Range("H136:N136").Value = $arr

$wsCUL = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 661.6667
$arr[0,1] = 661.6667
$arr[0,2] = 0
$arr[0,3] = 1985.0001
$arr[0,4] = 0
$arr[0,5] = -1445.0001
$arr[0,6] = $null
$wsCUL.Range("H59:N59").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 25000
$arr[0,1] = 25000
$arr[0,2] = 0
$arr[0,3] = 75000
$arr[0,4] = 0
$arr[0,5] = -73752
$arr[0,6] = $null
$wsCUL.Range("H87:N87").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 25000
$arr[0,1] = 25000
$arr[0,2] = 0
$arr[0,3] = 225000
$arr[0,4] = 0
$arr[0,5] = -218760
$arr[0,6] = $null
$wsCUL.Range("H90:N90").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1110.2
$arr[0,1] = 500
$arr[0,2] = 1262.75
$arr[0,3] = 1500
$arr[0,4] = 3788.25
$arr[0,5] = -621
$arr[0,6] = -5546.25
$wsCUL.Range("H103:N103").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = $null
$arr[0,6] = $null
$wsCUL.Range("H120:N120").Value = $arr

$wsLTW = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = $null
$arr[0,6] = $null
$wsLTW.Range("H7:N7").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 5128.4287
$arr[0,1] = 3699.889
$arr[0,2] = 6199.8335
$arr[0,3] = 3699.889
$arr[0,4] = 6199.8335
$arr[0,5] = -3511.889
$arr[0,6] = -6575.8335
$wsLTW.Range("H46:N46").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1514
$arr[0,1] = 1979.5
$arr[0,2] = 1048.5
$arr[0,3] = 1979.5
$arr[0,4] = 1048.5
$arr[0,5] = -731.5
$arr[0,6] = -3544.5
$wsLTW.Range("H93:N93").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = $null
$arr[0,6] = $null
$wsLTW.Range("H126:N126").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2090.1428
$arr[0,1] = 1881
$arr[0,2] = 2247
$arr[0,3] = 5643
$arr[0,4] = 6741
$arr[0,5] = -3093
$arr[0,6] = -11841
$wsLTW.Range("H136:N136").Value = $arr

$wsWVR = $wb.Worksheets.Item("WVR")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1544.6
$arr[0,1] = 1544.6
$arr[0,2] = 0
$arr[0,3] = 4633.799999999999
$arr[0,4] = 0
$arr[0,5] = -2103.799999999999
$arr[0,6] = $null
$wsWVR.Range("H132:N132").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2576.5
$arr[0,1] = 1810
$arr[0,2] = 5898
$arr[0,3] = 5430
$arr[0,4] = 17694
$arr[0,5] = -2880
$arr[0,6] = -22794
$wsWVR.Range("H136:N136").Value = $arr
